$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("D3").Value = "-"
$ws.Range("E3").Value = "-"
$ws.Range("F3").Value = "[Aline S. M.-T. M. Metalicos, Aline S. M.-T. M. Metalicos]"

# Row 4
$ws.Range("D4").Value = "[Ernane-Desenho tecnico mecanico, -]"
$ws.Range("F4").Value = "[Weslei-Metrologia 1, -, -, Emerson-Comandos Eletricos]"

# Row 6
$ws.Range("D6").Value = "[Ernane-Desenho tecnico mecanico, -]"
$ws.Range("E6").Value = "Anselmo-Gestao Intregr"
$ws.Range("F6").Value = "[Weslei-Metrologia 1, -, -, Emerson-Comandos Eletricos]"

# Row 7
$ws.Range("C7").Value = "-"
$ws.Range("D7").Value = "[Ernane-Desenho tecnico mecanico, -]"
$ws.Range("E7").Value = "Anselmo-Gestao Intregr"
$ws.Range("F7").Value = "[Weslei-Metrologia 1, -, -, Emerson-Comandos Eletricos]"

# Row 8
$ws.Range("F8").Value = "[Weslei-Metrologia 1, -, -, Emerson-Comandos Eletricos]"
